$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to remain plain text, even when
# it looks like a number (e.g. "0.519", "6.32"), without leaving a
# residual cell-level number-format style behind.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.664.53"
$ws.Range("E2").Value = "  -1.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.771.35"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "594.87"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "166.88"
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.768.90"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - XRP
Set-TextValue $ws.Range("D9") "0.519"
$ws.Range("E9").Value = "  -0.25%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.28%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "6.32"
$ws.Range("E11").Value = "  -1.86%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.45%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000254"
$ws.Range("E13").Value = "  -2.51%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "35.98"
$ws.Range("E14").Value = "  -0.39%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.402.40"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.756.31"
$ws.Range("E16").Value = "  +0.20%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.591.25"
$ws.Range("E17").Value = "  -1.44%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +1.76%  "

# Row 19 - now Polkadot (was TRON)
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D19") "6.98"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20 - now TRON (was Polkadot)
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D20") "0.111"
$ws.Range("E20").Value = "  -0.78%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "10.02"
$ws.Range("E21").Value = "  -6.98%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "456.40"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -0.50%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +4.54%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "83.13"
$ws.Range("E25").Value = "  -1.26%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D26") "11.89"
$ws.Range("E26").Value = "  -0.45%  "

# Row 27 - Fetch.AI
Set-TextValue $ws.Range("D27") "2.12"
$ws.Range("E27").Value = "  -2.53%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "10.06"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.14%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.34%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +2.17%  "

# Row 32 - NEARProtocol
Set-TextValue $ws.Range("D32") "7.24"
$ws.Range("E32").Value = "  -0.71%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "29.63"
$ws.Range("E33").Value = "  -1.12%  "

# Row 34 - Aptos
$ws.Range("E34").Value = "  -0.39%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("E35").Value = "  +0.12%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "3.723.33"
$ws.Range("E36").Value = "  +0.33%  "

# Row 37 - Hedera
Set-TextValue $ws.Range("D37") "0.100"
$ws.Range("E37").Value = "  -0.48%  "

# Row 38 - dogwifhat
Set-TextValue $ws.Range("D38") "3.30"
$ws.Range("E38").Value = "  -2.53%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.07%  "

# Row 40 - Mantle
Set-TextValue $ws.Range("D40") "0.994"
$ws.Range("E40").Value = "  -1.19%  "

# Row 41 - Filecoin
Set-TextValue $ws.Range("D41") "5.74"
$ws.Range("E41").Value = "  -0.75%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  -0.17%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  +0.00%  "

# Row 44 - Arweave
Set-TextValue $ws.Range("D44") "45.84"
$ws.Range("E44").Value = "  +5.91%  "

# Row 45 - OKB
Set-TextValue $ws.Range("D45") "48.18"
$ws.Range("E45").Value = "  +3.81%  "

# Row 46 - TheGraph
$ws.Range("E46").Value = "  -1.57%  "

# Row 47 - now Monero (was Cosmos)
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D47") "148.94"
$ws.Range("E47").Value = "  +2.18%  "

# Row 48 - now Cosmos (was Monero)
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "8.31"
$ws.Range("E48").Value = "  -2.17%  "

# Row 49 - Bittensor
Set-TextValue $ws.Range("D49") "388.10"
$ws.Range("E49").Value = "  -0.17%  "

# Row 50 - Stacks
Set-TextValue $ws.Range("D50") "1.81"
$ws.Range("E50").Value = "  -5.23%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "26.03"
$ws.Range("E51").Value = "  -0.25%  "
